$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 365 (shifts existing rows 365-383 down to 366-384)
$ws.Rows(365).Insert()

# Copy the static / repeated columns from the row above (row 364, now identical pattern)
$ws.Range("A365").Value = 3
$ws.Range("B365").Value = "Femacal de La Calera"
$ws.Range("C365").Value = "Coquimbo"
$ws.Range("D365").Value = 44753
$ws.Range("E365").Value = 5
$ws.Range("F365").Value = 100112040
$ws.Range("G365").Value = "Cilantro"
$ws.Range("H365").Value = "Sin especificar"
$ws.Range("I365").Value = "Primera"
$ws.Range("J365").Value = 310
$ws.Range("K365").Value = 3300
$ws.Range("L365").Value = 3500
$ws.Range("M365").Value = 3397
$ws.Range("N365").Value = "$/docena de atados (3 kilos)"
$ws.Range("O365").Value = "Provincia de Quillota"
$ws.Range("P365").Value = 1132
$ws.Range("Q365").Value = 3
$ws.Range("R365").Value = "Hortaliza"
